$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModelRuns")

# ------------------------------------------------------------------
# 1) Insert a new row at 80 (this copies formatting from row 79 above,
#    and pushes the old row 80 -> row 81, exactly like the diff shows)
# ------------------------------------------------------------------
$ws.Rows.Item(80).Insert()

# ------------------------------------------------------------------
# 2) Fill in the new row 80 with the "2035_TM160_IPA_13" run data
# ------------------------------------------------------------------
$ws.Range("A80").Value = 2035
$ws.Range("B80").Value = "2035_TM160_IPA_13"
$ws.Range("C80").Value = "RTP2025_IP"
$ws.Range("D80").Value = "IPA"
$ws.Range("E80").Value = "IPA with lower AOC"
$ws.Range("F80").Value = "FBP scaled to RGF"
$ws.Range("G80").Value = "run182"
$ws.Range("H80").Value = "current"
$ws.Range("I80").Value = "M:\Application\Model One\RTP2021\Blueprint\INPUT_DEVELOPMENT\Networks\BlueprintNetworks_64\net_2035_Blueprint_tollscsv"
$ws.Range("J80").Value = "model3-c"
$ws.Range("K80").Value = "https://app.asana.com/0/1204085012544660/1206153405312420/f"
$ws.Range("L80").Value = 13.68
$ws.Range("M80").Value = "na"
$ws.Range("N80").Value = "na"
$ws.Range("O80").Value = 0.87
$ws.Range("P80").Value = 0.78
$ws.Range("Q80").Value = 100
$ws.Range("R80").Value = 0
$ws.Range("S80").Value = 75

# Row insert carried the Hyperlink look-alike style into K80 (copied from
# K79 above which was a real hyperlink cell). Re-stamp K80's format from a
# normal (non-hyperlink) cell in the same row so it matches the rest of row 80.
$ws.Range("A80").Copy() | Out-Null
$ws.Range("K80").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# L80 carries the "highlighted" number style used for the most recent/changed
# AOC value in this table (same style as used e.g. in L71).
$ws.Range("L71").Copy() | Out-Null
$ws.Range("L80").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("L80").Value = 13.68

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) K79 no longer is a hyperlink -- remove the hyperlink and restore its
#    formatting to match the rest of row 79 (non-hyperlink look).
# ------------------------------------------------------------------
$ws.Range("K79").Hyperlinks.Delete()
$ws.Range("L79").Copy() | Out-Null
$ws.Range("K79").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# The "Hyperlink" cell style is no longer used anywhere in the workbook now,
# so remove it (matches removal of the Hyperlink font/cellStyle/cellStyleXf).
$wb.Styles.Item("Hyperlink").Delete()

# ------------------------------------------------------------------
# 4) Add the two new helper rows (83 & 84) with formulas referencing L80/L79
# ------------------------------------------------------------------
$ws.Range("L83").Formula = "=L80*2"
$ws.Range("M83").Formula = "=L79*2"

$ws.Range("L84").Formula = "=L80*3"
$ws.Range("M84").Formula = "=L79*3"

$wb.Save()
